$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These D-column cells hold purely numeric-looking price strings (e.g. "24.08").
# Force them to Text format first so Excel does not auto-convert the assigned
# string into a floating-point Value (matching the source data which stores
# these as plain text, not numbers).
$textCells = @("D5", "D8", "D11", "D14", "D17", "D18", "D20", "D22", "D23", "D24", "D25", "D27", "D31", "D32", "D33", "D36", "D39", "D40", "D41", "D45", "D46", "D48", "D49")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.225.39"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "1.590.35"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").Value = "213.44"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("D8").Value = "24.08"
$ws.Range("E8").Value = "  +7.29%  "
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").Value = "0.0890"
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").Value = "1.818.34"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "1.593.62"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").Value = "0.531"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").Value = "28.268.14"
$ws.Range("E16").Value = "  +2.98%  "
$ws.Range("D17").Value = "63.22"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").Value = "227.42"
$ws.Range("E18").Value = "  +2.70%  "
$ws.Range("D19").Value = "0.0₃0709"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").Value = "7.48"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").Value = "4.10"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").Value = "9.32"
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("D24").Value = "1.94"
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("D25").Value = "151.85"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").Value = "0.107"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  -1.57%  "
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").Value = "0.0473"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "3.23"
$ws.Range("E32").Value = "  -0.83%  "
$ws.Range("D33").Value = "3.15"
$ws.Range("E33").Value = "  -1.43%  "
$ws.Range("D34").Value = "1.399.44"
$ws.Range("E34").Value = "  -4.74%  "
$ws.Range("E35").Value = "  -2.74%  "
$ws.Range("D36").Value = "1.03"
$ws.Range("E36").Value = "  -8.32%  "
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("E38").Value = "  -0.35%  "
$ws.Range("D39").Value = "2.56"
$ws.Range("E39").Value = "  +5.34%  "
$ws.Range("D40").Value = "0.542"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "0.814"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("E43").Value = "  +4.90%  "
$ws.Range("E44").Value = "  -3.88%  "
$ws.Range("D45").Value = "0.982"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").Value = "64.30"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").Value = "1.729.49"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("D48").Value = "87.40"
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").Value = "2.13"
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("E51").Value = "  -0.56%  "
